# Generate Report for Handback
# This script updates the localization-status workbook to reflect a failed
# handback transform: the overall status text changes from
# "Ready for handoff" to "Handback transform failed", and an explanatory
# error message is written into the "Error Detail" column (K) of the last
# data row on both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$zhError = "Handback file name: hv5v3yvu.ktp is different with handoff file name: 2b08ac33-0475-49f0-b45a-962f8a653e59.64ddc41c2d54c9f4b3f05612a7a9e3f325d28f4b.zh-cn."
$deError  = "Handback file name: hv5v3yvu.ktp is different with handoff file name: 2b08ac33-0475-49f0-b45a-962f8a653e59.64ddc41c2d54c9f4b3f05612a7a9e3f325d28f4b.de-de."

# --- Overview sheet: update status for both locales on the last row ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B7").Value = $newStatus
$wsOverview.Range("C7").Value = $newStatus

# --- zh-cn sheet: update status and add the error detail ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C7").Value = $newStatus
$wsZhCn.Range("K7").Value = $zhError

# --- de-de sheet: update status and add the error detail ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C7").Value = $newStatus
$wsDeDe.Range("K7").Value = $deError
